# Update contact list: fix emails/messages, delete the last (4th) record,
# and widen the "message" column now that it holds a longer text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (id 1, ABOBAKAR): corrected email + message ---
$ws.Range("C2").Value = "ranaabobakarit@gmail.com"
$ws.Range("D2").Value = "kya hal ha"

# --- Row 3 (id 2, SHAHZEB): email swapped in, message updated ---
$ws.Range("C3").Value = "ranaabobakar777@gmail.com"
$ws.Range("D3").Value = "kya hal ha"

# --- Row 4 (id 3, ZEESHAN): new email, message updated ---
$ws.Range("C4").Value = "abobakarit786@gmail.com"
$ws.Range("D4").Value = "kya hal ha"

# The "message" cells on rows 3 & 4 get a distinct (plain black) font,
# separate from the default style used on row 2.
$ws.Range("D3:D4").Font.Color = 0

# --- Remove the 4th record (row 5) entirely ---
$ws.Range("A5:D5").EntireRow.Delete()

# Widen column D now that it carries the longer "message" text.
$ws.Columns.Item(4).ColumnWidth = 34.6

# Restore the selection to where the user left off.
[void]$ws.Range("D8").Select()
